$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the obsolete data rows (rows 3-5); only one data row remains afterwards.
$ws.Range("A3:T5").Delete() | Out-Null

# Update the remaining data row with the recalculated TPM-based values.
$ws.Range("A2").Value = "Resolving-Mac"
$ws.Range("B2").Value = "Ccl4"
$ws.Range("C2").Value = "Ccr3"
$ws.Range("D2").Value = "Resolving-Mac"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 66.00836066666666
$ws.Range("H2").Value = 198.025082
$ws.Range("I2").Value = 1
$ws.Range("J2").Value = 1
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 0.6666666666666666
$ws.Range("M2").Value = 0.1790523333333333
$ws.Range("N2").Value = 0.537157
$ws.Range("O2").Value = 1
$ws.Range("P2").Value = 1
$ws.Range("Q2").Value = 11.81895099687489
$ws.Range("R2").Value = 106.370558971874
$ws.Range("S2").Value = 1
$ws.Range("T2").Value = 1
